$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 ("I0") and J1 ("IF") - match the existing bold/
# centered/bordered header style used by the other header cells (row 1)
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-67: column I = I0, column J = IF
$data = @(@(3,6), @(6,6), @(6,7), @(9,9), @(5,6), @(4,5), @(8,8), @(7,8), @(7,8), @(8,8), @(9,9), @(7,8), @(9,9), @(9,9), @(6,6), @(1,2), @(8,8), @(6,7), @(7,7), @(5,6), @(6,7), @(7,7), @(6,7), @(8,8), @(5,6), @(7,7), @(5,5), @(9,9), @(5,6), @(5,6), @(8,9), @(8,8), @(7,7), @(5,6), @(6,7), @(4,5), @(6,7), @(6,7), @(10,10), @(6,7), @(6,7), @(7,7), @(6,6), @(7,7), @(8,9), @(8,8), @(6,6), @(9,9), @(9,9), @(5,5), @(5,6), @(6,7), @(6,9), @(4,5), @(5,7), @(9,9), @(7,7), @(5,5), @(6,7), @(7,8), @(8,9), @(4,4), @(9,9), @(6,6), @(5,5), @(6,6))

for ($k = 0; $k -lt $data.Count; $k++) {
    $r = $k + 2
    $ws.Cells.Item($r, 9).Value = $data[$k][0]
    $ws.Cells.Item($r, 10).Value = $data[$k][1]
}
